# 自动更新Excel文件 - 2026-02-17 23:21:51
# Decrement the "剩余" (remaining) count in column E by 1 for every data
# row (rows 2-99), except row 36 whose value stays the same.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) { continue }
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $cell.Value2 - 1
}
